$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.007.12'
Set-TextValue 'D3' '1.822.66'
Set-TextValue 'E3' '  -1.10%  '
Set-TextValue 'E4' '  -0.58%  '
Set-TextValue 'D5' '309.71'
Set-TextValue 'E6' '  -0.46%  '
Set-TextValue 'D7' '0.4618'
Set-TextValue 'E7' '  -3.01%  '
Set-TextValue 'D8' '0.3641'
Set-TextValue 'E8' '  -1.93%  '
Set-TextValue 'D9' '0.07284'
Set-TextValue 'E9' '  -2.48%  '
Set-TextValue 'D10' '0.8649'
Set-TextValue 'E10' '  -2.72%  '
Set-TextValue 'D12' '1.886.49'
Set-TextValue 'E12' '  +1.96%  '
Set-TextValue 'D13' '0.07606'
Set-TextValue 'E13' '  +2.95%  '
Set-TextValue 'D14' '93.23'
Set-TextValue 'E14' '  -0.30%  '
Set-TextValue 'D15' '5.330'
Set-TextValue 'E15' '  -2.83%  '
Set-TextValue 'D16' '6.494'
Set-TextValue 'E16' '  -1.60%  '
Set-TextValue 'E17' '  -0.72%  '
Set-TextValue 'D18' '0.000008626'
Set-TextValue 'E18' '  -2.51%  '
Set-TextValue 'D19' '1.008'
Set-TextValue 'E19' '  -0.54%  '
Set-TextValue 'D20' '27.420.64'
Set-TextValue 'E20' '  -0.15%  '
Set-TextValue 'E21' '  -2.40%  '
Set-TextValue 'E22' '  -3.55%  '
Set-TextValue 'E23' '  -1.52%  '
Set-TextValue 'D24' '2.113.37'
Set-TextValue 'E24' '  +1.81%  '
Set-TextValue 'D25' '151.63'
Set-TextValue 'E25' '  -0.63%  '
Set-TextValue 'D26' '1.856'
Set-TextValue 'E26' '  -2.50%  '
Set-TextValue 'D27' '18.22'
Set-TextValue 'E27' '  -2.39%  '
Set-TextValue 'D28' '2.090'
Set-TextValue 'E28' '  -3.56%  '
Set-TextValue 'D29' '5.100'
Set-TextValue 'E29' '  -3.61%  '
Set-TextValue 'D30' '115.90'
Set-TextValue 'E30' '  -1.96%  '
Set-TextValue 'D31' '0.08901'
Set-TextValue 'E31' '  -0.68%  '
Set-TextValue 'D32' '2.950'
Set-TextValue 'E32' '  -0.14%  '
Set-TextValue 'B33' 'ImmutableX'
Set-TextValue 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D33' '0.7277'
Set-TextValue 'E33' '  -4.28%  '
Set-TextValue 'B34' 'ARBITRUM'
Set-TextValue 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D34' '1.141'
Set-TextValue 'E34' '  -3.43%  '
Set-TextValue 'E35' '  -3.17%  '
Set-TextValue 'E36' '  -0.48%  '
Set-TextValue 'D37' '2.502'
Set-TextValue 'E37' '  +4.98%  '
Set-TextValue 'D38' '0.05274'
Set-TextValue 'E38' '  -1.66%  '
Set-TextValue 'E39' '  -3.11%  '
Set-TextValue 'E40' '  -2.50%  '
Set-TextValue 'D41' '2.925'
Set-TextValue 'E41' '  -2.53%  '
Set-TextValue 'D42' '7.171'
Set-TextValue 'E42' '  -1.99%  '
Set-TextValue 'E43' '  -2.95%  '
Set-TextValue 'D44' '0.1634'
Set-TextValue 'E44' '  -1.99%  '
Set-TextValue 'D45' '8.257'
Set-TextValue 'E45' '  -3.61%  '
Set-TextValue 'D46' '0.4858'
Set-TextValue 'E46' '  -2.54%  '
Set-TextValue 'E47' '  -0.54%  '
Set-TextValue 'D48' '10.11'
Set-TextValue 'E48' '  -4.89%  '
Set-TextValue 'D49' '103.16'
Set-TextValue 'E49' '  -1.99%  '
Set-TextValue 'D50' '1.632'
Set-TextValue 'E50' '  -3.26%  '
Set-TextValue 'D51' '0.06223'
Set-TextValue 'E51' '  -1.61%  '
